# Weekly refresh of the Jengibre (Hortaliza) price rows: each data row (2-20)
# is updated in place with the week's Fecha/Volumen/Precio values, matching
# the re-sorted weekly consolidation pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, Fecha(D, serial), Volumen(J), Precio min(K), Precio max(L), Precio prom(M), Precio por Kg(P)
$rows = @(
    @{ Row=2; D=44868; J=30; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=3; D=44930; J=30; K=17000; L=17000; M=17000; P=1308 },
    @{ Row=4; D=44841; J=30; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=5; D=44846; J=30; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=6; D=44922; J=30; K=17000; L=17000; M=17000; P=1308 },
    @{ Row=7; D=44810; J=40; K=12000; L=13000; M=12500; P=962 },
    @{ Row=8; D=44959; J=30; K=19000; L=19000; M=19000; P=1462 },
    @{ Row=9; D=44804; J=40; K=12000; L=13000; M=12500; P=962 },
    @{ Row=10; D=44874; J=30; K=17000; L=17000; M=17000; P=1308 },
    @{ Row=11; D=44797; J=60; K=12000; L=13000; M=12500; P=962 },
    @{ Row=12; D=44839; J=40; K=15000; L=16000; M=15500; P=1192 },
    @{ Row=13; D=44915; J=50; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=14; D=44859; J=30; K=13000; L=13000; M=13000; P=1000 },
    @{ Row=15; D=45155; J=30; K=20000; L=20000; M=20000; P=1538 },
    @{ Row=16; D=44943; J=30; K=17000; L=17000; M=17000; P=1308 },
    @{ Row=17; D=44832; J=60; K=17000; L=18000; M=17500; P=1346 },
    @{ Row=18; D=44880; J=30; K=17000; L=17000; M=17000; P=1308 },
    @{ Row=19; D=44894; J=30; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=20; D=44895; J=30; K=18000; L=18000; M=18000; P=1385 }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 4).Value = $row.D   # D: Fecha
    $ws.Cells.Item($row.Row, 10).Value = $row.J  # J: Volumen
    $ws.Cells.Item($row.Row, 11).Value = $row.K  # K: Precio minimo
    $ws.Cells.Item($row.Row, 12).Value = $row.L  # L: Precio maximo
    $ws.Cells.Item($row.Row, 13).Value = $row.M  # M: Precio promedio ponderado
    $ws.Cells.Item($row.Row, 16).Value = $row.P  # P: Precio por Kg
}
